# Update country stats (COVID-19 dataset) and the "last updated" banner.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Banner timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 22:22"

# Each entry: row, country name, Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes
$updates = @(
    @{ Row = 4; Pais = "Estados Unidos"; B = 391475; C = 24471; D = 21488; E = 357435; F = 9169; G = 1681; H = 12552 }
    @{ Row = 5; Pais = "España"; B = 140617; C = 3942; D = 43208; E = 83497; F = 7069; G = 571; H = 13912 }
    @{ Row = 16; Pais = "Canada"; B = 17847; C = 1180; D = 3935; E = 13537; F = 426; G = 52; H = 375 }
    @{ Row = 17; Pais = "Brasil"; B = 13717; C = 1534; D = 127; E = 12923; F = 296; G = 103; H = 667 }
    @{ Row = 18; Pais = "Austria"; B = 12635; C = 338; D = 4046; E = 8346; F = 243; G = 23; H = 243 }
    @{ Row = 97; Pais = "Costa de Marfil"; B = 349; C = 26; D = 41; E = 305; F = 0; G = 0; H = 3 }
    @{ Row = 101; Pais = "Niger"; B = 278; C = 25; D = 26; E = 241; F = 0; G = 1; H = 11 }
    @{ Row = 102; Pais = "San Marino"; B = 277; C = 0; D = 35; E = 210; F = 14; G = 0; H = 32 }
    @{ Row = 103; Pais = "Mauricio"; B = 268; C = 24; D = 8; E = 253; F = 3; G = 0; H = 7 }
    @{ Row = 104; Pais = "Estado de Palestina"; B = 261; C = 7; D = 42; E = 218; F = 0; G = 0; H = 1 }
    @{ Row = 121; Pais = "Guinea"; B = 144; C = 16; D = 5; E = 139; F = 0; G = 0; H = 0 }
    @{ Row = 122; Pais = "Guadalupe"; B = 139; C = 0; D = 31; E = 101; F = 14; G = 0; H = 7 }
    @{ Row = 123; Pais = "Brunei"; B = 135; C = 0; D = 85; E = 49; F = 3; G = 0; H = 1 }
    @{ Row = 127; Pais = "Trinidad yTobago"; B = 107; C = 2; D = 1; E = 98; F = 0; G = 0; H = 8 }
    @{ Row = 144; Pais = "Islas Caimanes"; B = 45; C = 6; D = 1; E = 43; F = 0; G = 0; H = 1 }
    @{ Row = 145; Pais = "Congo"; B = 45; C = 0; D = 2; E = 38; F = 0; G = 0; H = 5 }
    @{ Row = 146; Pais = "Macao"; B = 44; C = 0; D = 10; E = 34; F = 1; G = 0; H = 0 }
    @{ Row = 147; Pais = "San Martin (Parte Holandesa)"; B = 40; C = 3; D = 1; E = 33; F = 2; G = 0; H = 6 }
    @{ Row = 162; Pais = "Libia"; B = 20; C = 1; D = 1; E = 18; F = 0; G = 0; H = 1 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("A$r").Value = $u.Pais
    $ws.Range("B$r").Value = $u.B
    $ws.Range("C$r").Value = $u.C
    $ws.Range("D$r").Value = $u.D
    $ws.Range("E$r").Value = $u.E
    $ws.Range("F$r").Value = $u.F
    $ws.Range("G$r").Value = $u.G
    $ws.Range("H$r").Value = $u.H
}
